$d = $word.ActiveDocument

# Locate the empty paragraph that immediately follows the "Consumers:
# Metamask & Mist Browser" paragraph and immediately precedes the next
# empty paragraph, then delete it (paragraph mark included) so the two
# empty paragraphs collapse into one, matching the authored diff.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $i -gt 1) {
        $prev = $d.Paragraphs.Item($i - 1).Range.Text
        if ($prev -match "Mist Browser") {
            $p.Range.Delete()
            break
        }
    }
}
